$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2954.818
$ws.Range("I132").Value = 3100.3
$ws.Range("K132").Value = 9300.900000000001
$ws.Range("M132").Value = -6770.900000000001
$ws.Range("H137").Value = 2444.125
$ws.Range("J137").Value = 4001
$ws.Range("L137").Value = 12003
$ws.Range("N137").Value = -17103

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1617.7858
$ws.Range("I2").Value = 1005.8571
$ws.Range("K2").Value = 1005.8571
$ws.Range("M2").Value = -892.8570999999999
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H61").Value = 3469.2856
$ws.Range("I61").Value = 3469.2856
$ws.Range("K61").Value = 3469.2856
$ws.Range("M61").Value = -3257.2856
$ws.Range("H116").Value = 1617.7858
$ws.Range("I116").Value = 1005.8571
$ws.Range("K116").Value = 1005.8571
$ws.Range("M116").Value = 1288.1429
$ws.Range("H136").Value = 3469.2856
$ws.Range("I136").Value = 3469.2856
$ws.Range("K136").Value = 10407.8568
$ws.Range("M136").Value = -7857.856800000001
$ws.Range("H139").Value = 49999
$ws.Range("J139").Value = 49999
$ws.Range("L139").Value = 49999
$ws.Range("N139").Value = -60279

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1617.7858
$ws.Range("I3").Value = 1005.8571
$ws.Range("K3").Value = 1005.8571
$ws.Range("M3").Value = -891.8570999999999
$ws.Range("H31").Value = 500
$ws.Range("I31").Value = 500
$ws.Range("K31").Value = 500
$ws.Range("M31").Value = -248
$ws.Range("H33").Value = 999.5
$ws.Range("I33").Value = 999.5
$ws.Range("K33").Value = 999.5
$ws.Range("M33").Value = -663.5
$ws.Range("H36").Value = 5479
$ws.Range("I36").Value = 5479
$ws.Range("K36").Value = 5479
$ws.Range("M36").Value = -4945
$ws.Range("H37").Value = 1231
$ws.Range("I37").Value = 396.5
$ws.Range("J37").Value = 2900
$ws.Range("K37").Value = 396.5
$ws.Range("L37").Value = 2900
$ws.Range("M37").Value = -259.5
$ws.Range("N37").Value = -3174
$ws.Range("H46").Value = 2000
$ws.Range("J46").Value = 2000
$ws.Range("L46").Value = 2000
$ws.Range("N46").Value = -2596
$ws.Range("H51").Value = 100000
$ws.Range("J51").Value = 100000
$ws.Range("L51").Value = 100000
$ws.Range("N51").Value = -100982

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2373.1304
$ws.Range("I31").Value = 1914.9697
$ws.Range("K31").Value = 1914.9697
$ws.Range("M31").Value = -1619.9697
$ws.Range("H34").Value = 2373.1304
$ws.Range("I34").Value = 1914.9697
$ws.Range("K34").Value = 1914.9697
$ws.Range("M34").Value = -1712.9697
$ws.Range("H41").Value = 16011.8
$ws.Range("I41").Value = 7529.5
$ws.Range("J41").Value = 21666.666
$ws.Range("K41").Value = 7529.5
$ws.Range("L41").Value = 21666.666
$ws.Range("M41").Value = -7101.5
$ws.Range("N41").Value = -22522.666
$ws.Range("H50").Value = 27027.666
$ws.Range("H51").Value = 22525
$ws.Range("H58").Value = 6337.6665
$ws.Range("I58").Value = 4569
$ws.Range("K58").Value = 4569
$ws.Range("M58").Value = -4366
$ws.Range("H60").Value = 19500
$ws.Range("H61").Value = 22525
$ws.Range("H134").Value = 2873.1428
$ws.Range("I134").Value = 2873.1428
$ws.Range("K134").Value = 8619.428400000001
$ws.Range("M134").Value = -6084.428400000001
$ws.Range("H136").Value = 6337.6665
$ws.Range("I136").Value = 4569
$ws.Range("K136").Value = 13707
$ws.Range("M136").Value = -11157

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1933.4348
$ws.Range("I4").Value = 1909.9333
$ws.Range("J4").Value = 1977.5
$ws.Range("K4").Value = 5729.7999
$ws.Range("L4").Value = 5932.5
$ws.Range("M4").Value = -5617.7999
$ws.Range("N4").Value = -6156.5
$ws.Range("H14").Value = 25120.916
$ws.Range("I14").Value = 25120.916
$ws.Range("K14").Value = 75362.74800000001
$ws.Range("M14").Value = -75189.74800000001
$ws.Range("H64").Value = 1983.3334
$ws.Range("J64").Value = 2475
$ws.Range("L64").Value = 7425
$ws.Range("N64").Value = -7965
$ws.Range("H67").Value = 1983.3334
$ws.Range("J67").Value = 2475
$ws.Range("L67").Value = 7425
$ws.Range("N67").Value = -9297
$ws.Range("H68").Value = 1398.125
$ws.Range("J68").Value = 1237.25
$ws.Range("L68").Value = 3711.75
$ws.Range("N68").Value = -5333.75
$ws.Range("H71").Value = 1398.125
$ws.Range("J71").Value = 1237.25
$ws.Range("L71").Value = 11135.25
$ws.Range("N71").Value = -19247.25
$ws.Range("H112").Value = 745.6667
$ws.Range("I112").Value = 745.6667
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 2237.0001
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H114").Value = 181.5
$ws.Range("I114").Value = 163.5
$ws.Range("K114").Value = 490.5
$ws.Range("M114").Value = 2763.5
$ws.Range("H121").Value = 864.2
$ws.Range("I121").Value = 377.8
$ws.Range("J121").Value = 1107.4
$ws.Range("K121").Value = 1133.4
$ws.Range("L121").Value = 3322.2
$ws.Range("M121").Value = 176.5999999999999
$ws.Range("N121").Value = -5942.200000000001
$ws.Range("H123").Value = 2998.5
$ws.Range("I123").Value = 2997
$ws.Range("K123").Value = 8991
$ws.Range("M123").Value = -6541
$ws.Range("H138").Value = 6155.375
$ws.Range("J138").Value = 6713.25
$ws.Range("L138").Value = 20139.75
$ws.Range("N138").Value = -30419.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 1193.7778
$ws.Range("I5").Value = 1193.7778
$ws.Range("K5").Value = 1193.7778
$ws.Range("M5").Value = -1081.7778
$ws.Range("H122").Value = 3700
$ws.Range("I122").Value = 2320.2
$ws.Range("J122").Value = 5999.6665
$ws.Range("K122").Value = 6960.599999999999
$ws.Range("L122").Value = 17998.9995
$ws.Range("M122").Value = -4510.599999999999
$ws.Range("N122").Value = -22898.9995
$ws.Range("H132").Value = 2814.1428
$ws.Range("I132").Value = 2760.2
$ws.Range("K132").Value = 8280.599999999999
$ws.Range("M132").Value = -5750.599999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2602.3333
$ws.Range("I136").Value = 2602.3333
$ws.Range("K136").Value = 7806.999899999999
$ws.Range("M136").Value = -5256.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 749
$ws.Range("J107").Value = 999
$ws.Range("L107").Value = 2997
$ws.Range("N107").Value = -6837
$ws.Range("H119").Value = 16000
$ws.Range("J119").Value = 16000
$ws.Range("L119").Value = 16000
$ws.Range("N119").Value = -25676
$ws.Range("H136").Value = 10288.667
$ws.Range("I136").Value = 11299.667
$ws.Range("K136").Value = 33899.001
$ws.Range("M136").Value = -31349.001
